# Apply "Included mean in code for all experiments" changes.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Cartesian Velocity - Linear" -------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Tiny floating point precision refresh on existing Standard Deviation / Maximum rows
$ws1.Cells.Item(3,3).Value = 0.06852659971056689
$ws1.Cells.Item(3,4).Value = 0.07208002718491419
$ws1.Cells.Item(3,5).Value = 0.03338962121225973

$ws1.Cells.Item(4,2).Value = 0.5080253236132926
$ws1.Cells.Item(4,3).Value = 0.432311073253249
$ws1.Cells.Item(4,4).Value = 0.4130668033357447
$ws1.Cells.Item(4,5).Value = 0.3708904621868302

# New "Mean" row
$ws1.Cells.Item(5,1).Value = "Mean"
$ws1.Cells.Item(5,2).Value = 0.09519282655710233
$ws1.Cells.Item(5,3).Value = 0.09110913624555038
$ws1.Cells.Item(5,4).Value = 0.09406766988253523
$ws1.Cells.Item(5,5).Value = 0.013976850285178

# --- Sheet 2: "Vertical" --------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(4,2).Value = 0.3646477683989287

$ws2.Cells.Item(5,1).Value = "Mean"
$ws2.Cells.Item(5,2).Value = 0.01951392687395583

# --- Sheet 3: "Updated Lin" -----------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(5,1).Value = "Mean"
$ws3.Cells.Item(5,2).Value = 0.1431671287973773
